$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-01-30 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-01-31 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("741÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "306÷4=", 2) | Out-Null
$d.Content.Find.Execute("435÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "743÷6=", 2) | Out-Null
$d.Content.Find.Execute("894÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "360÷4=", 2) | Out-Null
$d.Content.Find.Execute("764÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "923÷2=", 2) | Out-Null
$d.Content.Find.Execute("928÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "473÷2=", 2) | Out-Null
$d.Content.Find.Execute("208÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "984÷4=", 2) | Out-Null
$d.Content.Find.Execute("373÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "278÷9=", 2) | Out-Null
$d.Content.Find.Execute("140÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "230÷3=", 2) | Out-Null
$d.Content.Find.Execute("363÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "429÷9=", 2) | Out-Null
$d.Content.Find.Execute("275÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "604÷2=", 2) | Out-Null
$d.Content.Find.Execute("713÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "335÷2=", 2) | Out-Null
$d.Content.Find.Execute("558÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "728÷3=", 2) | Out-Null
$d.Content.Find.Execute("921÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "197÷7=", 2) | Out-Null
$d.Content.Find.Execute("207÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "562÷9=", 2) | Out-Null
$d.Content.Find.Execute("751÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "124÷4=", 2) | Out-Null
$d.Content.Find.Execute("921÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "151÷4=", 2) | Out-Null
$d.Content.Find.Execute("465÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "301÷9=", 2) | Out-Null
$d.Content.Find.Execute("408÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "576÷3=", 2) | Out-Null
$d.Content.Find.Execute("554÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "817÷5=", 2) | Out-Null
$d.Content.Find.Execute("711÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "857÷6=", 2) | Out-Null
$d.Content.Find.Execute("681÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "948÷2=", 2) | Out-Null
$d.Content.Find.Execute("782÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "676÷4=", 2) | Out-Null
$d.Content.Find.Execute("244÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "556÷4=", 2) | Out-Null
$d.Content.Find.Execute("100÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "458÷6=", 2) | Out-Null
$d.Content.Find.Execute("529÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "593÷3=", 2) | Out-Null
